$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix C27: 1600 -> 1700
$ws.Cells.Item(27, 3).Value = 1700

# Append new data rows 28-58 (entered in the order they were authored)
$ws.Cells.Item(28, 1).Value = "Ahmet"
$ws.Cells.Item(28, 2).Value = 3
$ws.Cells.Item(28, 3).Value = 1900
$ws.Cells.Item(28, 4).Value = 3
$ws.Cells.Item(28, 5).Value = 2
$ws.Cells.Item(28, 6).Value = 7
$ws.Cells.Item(28, 7).Value = 4
$ws.Cells.Item(28, 8).Value = 13
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 10).Value = 4
$ws.Cells.Item(28, 11).Value = 1
$ws.Cells.Item(28, 12).Value = 8

$ws.Cells.Item(29, 1).Value = "Doğukan"
$ws.Cells.Item(29, 2).Value = 2.5
$ws.Cells.Item(29, 3).Value = 1800
$ws.Cells.Item(29, 4).Value = 3
$ws.Cells.Item(29, 5).Value = 4
$ws.Cells.Item(29, 6).Value = 5
$ws.Cells.Item(29, 7).Value = 8
$ws.Cells.Item(29, 8).Value = 7
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 8
$ws.Cells.Item(29, 11).Value = 1
$ws.Cells.Item(29, 12).Value = 5

$ws.Cells.Item(30, 1).Value = "Tuna"
$ws.Cells.Item(30, 2).Value = 2
$ws.Cells.Item(30, 3).Value = 1700
$ws.Cells.Item(30, 4).Value = 3
$ws.Cells.Item(30, 5).Value = 2
$ws.Cells.Item(30, 6).Value = 5
$ws.Cells.Item(30, 7).Value = 8
$ws.Cells.Item(30, 8).Value = 6
$ws.Cells.Item(30, 9).Value = 2
$ws.Cells.Item(30, 10).Value = 5
$ws.Cells.Item(30, 11).Value = 2
$ws.Cells.Item(30, 12).Value = 6

$ws.Cells.Item(31, 1).Value = "Meriç"
$ws.Cells.Item(31, 2).Value = 3
$ws.Cells.Item(31, 3).Value = 2100
$ws.Cells.Item(31, 4).Value = 3
$ws.Cells.Item(31, 5).Value = 4
$ws.Cells.Item(31, 6).Value = 6
$ws.Cells.Item(31, 7).Value = 7
$ws.Cells.Item(31, 8).Value = 18
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 5
$ws.Cells.Item(31, 11).Value = 2
$ws.Cells.Item(31, 12).Value = 6

$ws.Cells.Item(32, 1).Value = "Kadir"
$ws.Cells.Item(32, 2).Value = 2
$ws.Cells.Item(32, 3).Value = 1500
$ws.Cells.Item(32, 4).Value = 4
$ws.Cells.Item(32, 5).Value = 5
$ws.Cells.Item(32, 6).Value = 3
$ws.Cells.Item(32, 7).Value = 6
$ws.Cells.Item(32, 8).Value = 7.5
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 10).Value = 10
$ws.Cells.Item(32, 11).Value = 0
$ws.Cells.Item(32, 12).Value = 5

$ws.Cells.Item(33, 1).Value = "İbrahim"
$ws.Cells.Item(33, 2).Value = 2
$ws.Cells.Item(33, 3).Value = 2000
$ws.Cells.Item(33, 4).Value = 3
$ws.Cells.Item(33, 5).Value = 5
$ws.Cells.Item(33, 6).Value = 8
$ws.Cells.Item(33, 7).Value = 7
$ws.Cells.Item(33, 8).Value = 23
$ws.Cells.Item(33, 9).Value = 0.5
$ws.Cells.Item(33, 10).Value = 7
$ws.Cells.Item(33, 11).Value = 1
$ws.Cells.Item(33, 12).Value = 3

$ws.Cells.Item(34, 1).Value = "Ilgın"
$ws.Cells.Item(34, 2).Value = 2
$ws.Cells.Item(34, 3).Value = 2000
$ws.Cells.Item(34, 4).Value = 3
$ws.Cells.Item(34, 5).Value = 6
$ws.Cells.Item(34, 6).Value = 6
$ws.Cells.Item(34, 7).Value = 8
$ws.Cells.Item(34, 8).Value = 12
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 9
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 12).Value = 5

$ws.Cells.Item(35, 1).Value = "Çetin"
$ws.Cells.Item(35, 2).Value = 2
$ws.Cells.Item(35, 3).Value = 1500
$ws.Cells.Item(35, 4).Value = 0
$ws.Cells.Item(35, 5).Value = 3.5
$ws.Cells.Item(35, 6).Value = 3
$ws.Cells.Item(35, 7).Value = 5
$ws.Cells.Item(35, 8).Value = 10
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(35, 10).Value = 7
$ws.Cells.Item(35, 11).Value = 0
$ws.Cells.Item(35, 12).Value = 0

$ws.Cells.Item(36, 1).Value = "Efe"
$ws.Cells.Item(36, 2).Value = 2
$ws.Cells.Item(36, 3).Value = 2300
$ws.Cells.Item(36, 4).Value = 1
$ws.Cells.Item(36, 5).Value = 3
$ws.Cells.Item(36, 6).Value = 4
$ws.Cells.Item(36, 7).Value = 6
$ws.Cells.Item(36, 8).Value = 5
$ws.Cells.Item(36, 9).Value = 1
$ws.Cells.Item(36, 10).Value = 5
$ws.Cells.Item(36, 11).Value = 0
$ws.Cells.Item(36, 12).Value = 3

$ws.Cells.Item(37, 1).Value = "Buse"
$ws.Cells.Item(37, 2).Value = 2
$ws.Cells.Item(37, 3).Value = 1300
$ws.Cells.Item(37, 4).Value = 4
$ws.Cells.Item(37, 5).Value = 8
$ws.Cells.Item(37, 6).Value = 7
$ws.Cells.Item(37, 7).Value = 5
$ws.Cells.Item(37, 8).Value = 30
$ws.Cells.Item(37, 9).Value = 0
$ws.Cells.Item(37, 10).Value = 5
$ws.Cells.Item(37, 11).Value = 2
$ws.Cells.Item(37, 12).Value = 8

$ws.Cells.Item(38, 1).Value = "Kaan"
$ws.Cells.Item(38, 2).Value = 2.5
$ws.Cells.Item(38, 3).Value = 1700
$ws.Cells.Item(38, 4).Value = 3
$ws.Cells.Item(38, 5).Value = 6
$ws.Cells.Item(38, 6).Value = 5
$ws.Cells.Item(38, 7).Value = 9
$ws.Cells.Item(38, 8).Value = 3
$ws.Cells.Item(38, 9).Value = 0
$ws.Cells.Item(38, 10).Value = 7
$ws.Cells.Item(38, 11).Value = 3
$ws.Cells.Item(38, 12).Value = 7

$ws.Cells.Item(39, 1).Value = "Çağdaş"
$ws.Cells.Item(39, 2).Value = 1
$ws.Cells.Item(39, 3).Value = 1900
$ws.Cells.Item(39, 4).Value = 4
$ws.Cells.Item(39, 5).Value = 6
$ws.Cells.Item(39, 6).Value = 6
$ws.Cells.Item(39, 7).Value = 5
$ws.Cells.Item(39, 8).Value = 22
$ws.Cells.Item(39, 9).Value = 0
$ws.Cells.Item(39, 10).Value = 5
$ws.Cells.Item(39, 11).Value = 2
$ws.Cells.Item(39, 12).Value = 4

$ws.Cells.Item(40, 1).Value = "Eshabil"
$ws.Cells.Item(40, 2).Value = 2.5
$ws.Cells.Item(40, 3).Value = 1800
$ws.Cells.Item(40, 4).Value = 5
$ws.Cells.Item(40, 5).Value = 6
$ws.Cells.Item(40, 6).Value = 6
$ws.Cells.Item(40, 7).Value = 8
$ws.Cells.Item(40, 8).Value = 42
$ws.Cells.Item(40, 9).Value = 1
$ws.Cells.Item(40, 10).Value = 7
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).Value = 6

$ws.Cells.Item(41, 1).Value = "Mehmet"
$ws.Cells.Item(41, 2).Value = 3
$ws.Cells.Item(41, 3).Value = 1500
$ws.Cells.Item(41, 4).Value = 1
$ws.Cells.Item(41, 5).Value = 3
$ws.Cells.Item(41, 6).Value = 5
$ws.Cells.Item(41, 7).Value = 6
$ws.Cells.Item(41, 8).Value = 14
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 10).Value = 4
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 12).Value = 6

$ws.Cells.Item(42, 1).Value = "Buğrahan"
$ws.Cells.Item(42, 2).Value = 3
$ws.Cells.Item(42, 3).Value = 1300
$ws.Cells.Item(42, 4).Value = 4
$ws.Cells.Item(42, 5).Value = 3
$ws.Cells.Item(42, 6).Value = 8
$ws.Cells.Item(42, 7).Value = 8
$ws.Cells.Item(42, 8).Value = 8
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 6
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 12).Value = 5

$ws.Cells.Item(43, 1).Value = "Mustafa"
$ws.Cells.Item(43, 2).Value = 2
$ws.Cells.Item(43, 3).Value = 2000
$ws.Cells.Item(43, 4).Value = 3
$ws.Cells.Item(43, 5).Value = 3
$ws.Cells.Item(43, 6).Value = 5
$ws.Cells.Item(43, 7).Value = 7
$ws.Cells.Item(43, 8).Value = 7
$ws.Cells.Item(43, 9).Value = 1
$ws.Cells.Item(43, 10).Value = 7
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 6

$ws.Cells.Item(44, 1).Value = "Emre"
$ws.Cells.Item(44, 2).Value = 3
$ws.Cells.Item(44, 3).Value = 1950
$ws.Cells.Item(44, 4).Value = 2
$ws.Cells.Item(44, 5).Value = 6
$ws.Cells.Item(44, 6).Value = 5
$ws.Cells.Item(44, 7).Value = 6
$ws.Cells.Item(44, 8).Value = 13
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 5
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 12).Value = 4

$ws.Cells.Item(45, 1).Value = "Yeşim"
$ws.Cells.Item(45, 2).Value = 2
$ws.Cells.Item(45, 3).Value = 1350
$ws.Cells.Item(45, 4).Value = 4
$ws.Cells.Item(45, 5).Value = 5
$ws.Cells.Item(45, 6).Value = 7
$ws.Cells.Item(45, 7).Value = 5
$ws.Cells.Item(45, 8).Value = 8
$ws.Cells.Item(45, 9).Value = 1
$ws.Cells.Item(45, 10).Value = 10
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 12).Value = 8

$ws.Cells.Item(46, 1).Value = "Mümin"
$ws.Cells.Item(46, 2).Value = 2
$ws.Cells.Item(46, 3).Value = 2400
$ws.Cells.Item(46, 4).Value = 3
$ws.Cells.Item(46, 5).Value = 6
$ws.Cells.Item(46, 6).Value = 5
$ws.Cells.Item(46, 7).Value = 7
$ws.Cells.Item(46, 8).Value = 20
$ws.Cells.Item(46, 9).Value = 1
$ws.Cells.Item(46, 10).Value = 8
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 2

$ws.Cells.Item(47, 1).Value = "Beril"
$ws.Cells.Item(47, 2).Value = 2
$ws.Cells.Item(47, 3).Value = 1500
$ws.Cells.Item(47, 4).Value = 4
$ws.Cells.Item(47, 5).Value = 6
$ws.Cells.Item(47, 6).Value = 3
$ws.Cells.Item(47, 7).Value = 10
$ws.Cells.Item(47, 8).Value = 6
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 10).Value = 10
$ws.Cells.Item(47, 11).Value = 2
$ws.Cells.Item(47, 12).Value = 8

$ws.Cells.Item(49, 1).Value = "Hatice"
$ws.Cells.Item(49, 2).Value = 2
$ws.Cells.Item(49, 3).Value = 1500
$ws.Cells.Item(49, 4).Value = 1
$ws.Cells.Item(49, 5).Value = 2
$ws.Cells.Item(49, 6).Value = 6
$ws.Cells.Item(49, 7).Value = 6
$ws.Cells.Item(49, 8).Value = 4
$ws.Cells.Item(49, 9).Value = 2
$ws.Cells.Item(49, 10).Value = 6
$ws.Cells.Item(49, 11).Value = 0
$ws.Cells.Item(49, 12).Value = 6

$ws.Cells.Item(50, 1).Value = "Selen"
$ws.Cells.Item(50, 2).Value = 2
$ws.Cells.Item(50, 3).Value = 1300
$ws.Cells.Item(50, 4).Value = 5
$ws.Cells.Item(50, 5).Value = 7
$ws.Cells.Item(50, 6).Value = 6
$ws.Cells.Item(50, 7).Value = 5
$ws.Cells.Item(50, 8).Value = 15
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 7
$ws.Cells.Item(50, 11).Value = 2
$ws.Cells.Item(50, 12).Value = 7

$ws.Cells.Item(51, 1).Value = "Esma"
$ws.Cells.Item(51, 2).Value = 2.5
$ws.Cells.Item(51, 3).Value = 1700
$ws.Cells.Item(51, 4).Value = 4
$ws.Cells.Item(51, 5).Value = 5.5
$ws.Cells.Item(51, 6).Value = 7
$ws.Cells.Item(51, 7).Value = 8
$ws.Cells.Item(51, 8).Value = 3
$ws.Cells.Item(51, 9).Value = 1
$ws.Cells.Item(51, 10).Value = 7
$ws.Cells.Item(51, 11).Value = 1
$ws.Cells.Item(51, 12).Value = 6

$ws.Cells.Item(52, 1).Value = "Cansın"
$ws.Cells.Item(52, 2).Value = 1
$ws.Cells.Item(52, 3).Value = 2200
$ws.Cells.Item(52, 4).Value = 6
$ws.Cells.Item(52, 5).Value = 5
$ws.Cells.Item(52, 6).Value = 6
$ws.Cells.Item(52, 7).Value = 4
$ws.Cells.Item(52, 8).Value = 17
$ws.Cells.Item(52, 9).Value = 1
$ws.Cells.Item(52, 10).Value = 5
$ws.Cells.Item(52, 11).Value = 0
$ws.Cells.Item(52, 12).Value = 0

$ws.Cells.Item(53, 1).Value = "Talha"
$ws.Cells.Item(53, 2).Value = 2.5
$ws.Cells.Item(53, 3).Value = 1600
$ws.Cells.Item(53, 4).Value = 5
$ws.Cells.Item(53, 5).Value = 6
$ws.Cells.Item(53, 6).Value = 5
$ws.Cells.Item(53, 7).Value = 3
$ws.Cells.Item(53, 8).Value = 30
$ws.Cells.Item(53, 9).Value = 5
$ws.Cells.Item(53, 10).Value = 5
$ws.Cells.Item(53, 11).Value = 2
$ws.Cells.Item(53, 12).Value = 6

$ws.Cells.Item(54, 1).Value = "Hasan"
$ws.Cells.Item(54, 2).Value = 3
$ws.Cells.Item(54, 3).Value = 1800
$ws.Cells.Item(54, 4).Value = 4
$ws.Cells.Item(54, 5).Value = 5
$ws.Cells.Item(54, 6).Value = 7
$ws.Cells.Item(54, 7).Value = 8
$ws.Cells.Item(54, 8).Value = 7
$ws.Cells.Item(54, 9).Value = 0
$ws.Cells.Item(54, 10).Value = 5
$ws.Cells.Item(54, 11).Value = 0
$ws.Cells.Item(54, 12).Value = 5

$ws.Cells.Item(55, 1).Value = "Harun"
$ws.Cells.Item(55, 2).Value = 2
$ws.Cells.Item(55, 3).Value = 2000
$ws.Cells.Item(55, 4).Value = 3
$ws.Cells.Item(55, 5).Value = 2
$ws.Cells.Item(55, 6).Value = 5
$ws.Cells.Item(55, 7).Value = 8
$ws.Cells.Item(55, 8).Value = 9
$ws.Cells.Item(55, 9).Value = 0.5
$ws.Cells.Item(55, 10).Value = 6
$ws.Cells.Item(55, 11).Value = 0
$ws.Cells.Item(55, 12).Value = 8

$ws.Cells.Item(56, 1).Value = "Serhat"
$ws.Cells.Item(56, 2).Value = 3
$ws.Cells.Item(56, 3).Value = 2300
$ws.Cells.Item(56, 4).Value = 2
$ws.Cells.Item(56, 5).Value = 6
$ws.Cells.Item(56, 6).Value = 3
$ws.Cells.Item(56, 7).Value = 7
$ws.Cells.Item(56, 8).Value = 45
$ws.Cells.Item(56, 9).Value = 0
$ws.Cells.Item(56, 10).Value = 6
$ws.Cells.Item(56, 11).Value = 1
$ws.Cells.Item(56, 12).Value = 4

$ws.Cells.Item(57, 1).Value = "Berk"
$ws.Cells.Item(57, 2).Value = 2
$ws.Cells.Item(57, 3).Value = 1500
$ws.Cells.Item(57, 4).Value = 2
$ws.Cells.Item(57, 5).Value = 5
$ws.Cells.Item(57, 6).Value = 5
$ws.Cells.Item(57, 7).Value = 7
$ws.Cells.Item(57, 8).Value = 8
$ws.Cells.Item(57, 9).Value = 0
$ws.Cells.Item(57, 10).Value = 9
$ws.Cells.Item(57, 11).Value = 0
$ws.Cells.Item(57, 12).Value = 7

$ws.Cells.Item(58, 1).Value = "Berat"
$ws.Cells.Item(58, 2).Value = 2
$ws.Cells.Item(58, 3).Value = 2300
$ws.Cells.Item(58, 4).Value = 2
$ws.Cells.Item(58, 5).Value = 7
$ws.Cells.Item(58, 6).Value = 6
$ws.Cells.Item(58, 7).Value = 7
$ws.Cells.Item(58, 8).Value = 21
$ws.Cells.Item(58, 9).Value = 0.5
$ws.Cells.Item(58, 10).Value = 7
$ws.Cells.Item(58, 11).Value = 1
$ws.Cells.Item(58, 12).Value = 5

$ws.Cells.Item(48, 1).Value = "Erdem"
$ws.Cells.Item(48, 2).Value = 2
$ws.Cells.Item(48, 3).Value = 1600
$ws.Cells.Item(48, 4).Value = 1
$ws.Cells.Item(48, 5).Value = 5
$ws.Cells.Item(48, 6).Value = 3
$ws.Cells.Item(48, 7).Value = 6.5
$ws.Cells.Item(48, 8).Value = 6
$ws.Cells.Item(48, 9).Value = 0
$ws.Cells.Item(48, 10).Value = 7
$ws.Cells.Item(48, 11).Value = 0
$ws.Cells.Item(48, 12).Value = 10

# Move selection to reflect the post-edit cursor position, as in the saved workbook
$ws.Range("L59").Select() | Out-Null

Write-Output "done"
